$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.05

$ws.Range("H3").Value = 2.75
$ws.Range("I3").Value = 2.7

$ws.Range("G5").Value = 1.83
$ws.Range("L5").Value = 1.45
$ws.Range("M5").Value = 2.37
$ws.Range("Q5").Value = 2.27
$ws.Range("X5").Value = 17.5
$ws.Range("Z5").Value = 6.7
$ws.Range("AA5").Value = 6.3
$ws.Range("AC5").Value = 150
$ws.Range("AE5").Value = 9

$ws.Range("T6").Value = 5.8
$ws.Range("W6").Value = 23
$ws.Range("AB6").Value = 19
$ws.Range("AE6").Value = 7.1
$ws.Range("AJ6").Value = 60
